# Re-style the three data tables (slides 14, 15 and 16) so they use the
# built-in "Medium Style 2 - Accent 1" table style instead of the custom
# Google-Slides-imported table style.
#
# Table.Style is a read-only reflection of the table's current style; the
# PowerPoint object model updates it via Table.ApplyStyle(StyleId).

$p = $ppt.ActivePresentation

$newStyleId = "{60246129-6BC9-40E4-847A-4D080DFDDACA}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)

    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}
